$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial number (45614 -> 2024-11-18).
# Update rows 2 through 33 to the next day (45615 -> 2024-11-19).
for ($row = 2; $row -le 33; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45614) {
        $cell.Value = 45615
    }
}
